$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Hebron"
$ws.Range("B5").Value = 42.3244
$ws.Range("C5").Value = -88.4524
$ws.Range("D5").Value = 1060
$ws.Range("D5").HorizontalAlignment = -4131

$ws.Range("A5:D5").Select()
